{"js": "// Merge the split \"<id>...</id>\" runs in each matching paragraph into a single\n// run, using the formatting of the first (\"<id>\") run, per the commit's\n// canonical OOXML diff (tc_p032v.docx): the <id> opening tag run, the middle\n// id-value run(s), and the closing </id> run get combined into one run\n// whose text is \"<id>VALUE</id>\".\n\nconst body = context.document.body;\nconst searchResults = body.search(\"<id>\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  const hit = searchResults.items[i];\n  const paragraph = hit.paragraphs.getFirst();\n  paragraph.load(\"text\");\n  await context.sync();\n\n  const fullText = paragraph.text;\n  // Only touch paragraphs that look like \"<id>...</id>\" (defensive check).\n  if (!/^<id>.*<\\/id>$/.test(fullText)) {\n    continue;\n  }\n\n  const paragraphRange = paragraph.getRange();\n  // Replacing the whole paragraph range's text collapses every run it spans\n  // into a single run. Word JS keeps the formatting of the first run in the\n  // replaced range, matching the <id> run's Courier New / 7f6000 / 18pt look.\n  paragraphRange.insertText(fullText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Merge the split \"<id>...</id>\" runs in each matching paragraph into a\n# single run, using the formatting of the first (\"<id>\") run, per the\n# commit's canonical OOXML diff (tc_p032v.docx): the <id> opening-tag run,\n# the middle id-value run(s), and the closing </id> run get combined into\n# one run whose text is \"<id>VALUE</id>\".\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"<id>\"\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n\nwhile ($rng.Find.Execute()) {\n    $para = $rng.Paragraphs(1)\n    $paraRange = $para.Range\n    # Exclude the trailing paragraph mark from the range we rewrite.\n    [void]$paraRange.MoveEnd(1, -1)\n    $fullText = $paraRange.Text\n\n    if ($fullText -like \"<id>*</id>\") {\n        # Re-assigning the exact same string is a no-op for the engine, so\n        # round-trip through a placeholder value first; this forces Word to\n        # collapse every run spanned by $paraRange into a single run that\n        # inherits the formatting of the first (\"<id>\") run, exactly as seen\n        # in the target OOXML.\n        $paraRange.Text = \"__MERGE_PLACEHOLDER__\"\n        $paraRange2 = $para.Range\n        [void]$paraRange2.MoveEnd(1, -1)\n        $paraRange2.Text = $fullText\n    }\n\n    # Move past this paragraph before searching again.\n    [void]$rng.Collapse(0)\n    [void]$rng.MoveEnd(1, 0)\n}\n"}
